$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Component" column (J): header in J1, then per-row stratigraphic
# component (lower / upper / sterile) for rows 2-32.
$ws.Range("J1").Value = "Component"
$ws.Range("J2").Value = "lower"
$ws.Range("J3").Value = "upper"
$ws.Range("J4").Value = "lower"
$ws.Range("J5").Value = "upper"
$ws.Range("J6").Value = "sterile"
$ws.Range("J7").Value = "sterile"
$ws.Range("J8").Value = "sterile"
$ws.Range("J9").Value = "upper"
$ws.Range("J10").Value = "upper"
$ws.Range("J11").Value = "upper"
$ws.Range("J12").Value = "upper"
$ws.Range("J13").Value = "upper"
$ws.Range("J14").Value = "upper"
$ws.Range("J15").Value = "upper"
$ws.Range("J16").Value = "upper"
$ws.Range("J17").Value = "upper"
$ws.Range("J18").Value = "lower"
$ws.Range("J19").Value = "lower"
$ws.Range("J20").Value = "upper"
$ws.Range("J21").Value = "upper"
$ws.Range("J22").Value = "lower"
$ws.Range("J23").Value = "upper"
$ws.Range("J24").Value = "upper"
$ws.Range("J25").Value = "upper"
$ws.Range("J26").Value = "lower"
$ws.Range("J27").Value = "lower"
$ws.Range("J28").Value = "upper"
$ws.Range("J29").Value = "lower"
$ws.Range("J30").Value = "lower"
$ws.Range("J31").Value = "sterile"
$ws.Range("J32").Value = "lower"

# Widen the new column to fit its contents, matching the authored layout.
$ws.Columns.Item(10).ColumnWidth = 11

# Scroll the view down and move the active selection, matching the
# author's final cursor position in the sheet.
$ws.Range("A11").Select()
$ws.Range("J16").Select()
